$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 303, pushing the existing data
# (old rows 303-442) down to rows 305-444.
$ws.Rows.Item(303).Insert()
$ws.Rows.Item(303).Insert()

# New row 303: Coliflor, "Primera" quality, week of 2021-11-10 (serial 44510)
$ws.Range("A303").Value = 8
$ws.Range("B303").Value = "Terminal La Palmera de La Serena"
$ws.Range("C303").Value = "Coquimbo"
$ws.Range("D303").Value = 44510
$ws.Range("E303").Value = 4
$ws.Range("F303").Value = 100112008
$ws.Range("G303").Value = "Coliflor"
$ws.Range("H303").Value = "Sin especificar"
$ws.Range("I303").Value = "Primera"
$ws.Range("J303").Value = 3300
$ws.Range("K303").Value = 600
$ws.Range("L303").Value = 700
$ws.Range("M303").Value = 650
$ws.Range("N303").Value = '$/unidad'
$ws.Range("O303").Value = "Provincia del Elquí"
$ws.Range("P303").Value = 650
$ws.Range("Q303").Value = 1
$ws.Range("R303").Value = "Hortaliza"

# New row 304: Coliflor, "Segunda" quality, week of 2021-11-10 (serial 44510)
$ws.Range("A304").Value = 8
$ws.Range("B304").Value = "Terminal La Palmera de La Serena"
$ws.Range("C304").Value = "Coquimbo"
$ws.Range("D304").Value = 44510
$ws.Range("E304").Value = 4
$ws.Range("F304").Value = 100112008
$ws.Range("G304").Value = "Coliflor"
$ws.Range("H304").Value = "Sin especificar"
$ws.Range("I304").Value = "Segunda"
$ws.Range("J304").Value = 1660
$ws.Range("K304").Value = 500
$ws.Range("L304").Value = 550
$ws.Range("M304").Value = 525
$ws.Range("N304").Value = '$/unidad'
$ws.Range("O304").Value = "Provincia del Elquí"
$ws.Range("P304").Value = 525
$ws.Range("Q304").Value = 1
$ws.Range("R304").Value = "Hortaliza"
